# Calculate BR values for all Analysis unit output.
#
# Add a new row to the "Library_Formula" sheet describing the new
# "getBRFieldName" library formula, following the same layout used by the
# existing rows: Action | Library | Formula Name | (Description is left
# blank) | Output | Input.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")

$row = 21

$ws.Cells.Item($row, 1).Value2 = "CREATE/MODIFY"
$ws.Cells.Item($row, 2).Value2 = "LIB_EWS"
$ws.Cells.Item($row, 3).Value2 = "getBRFieldName"
$ws.Cells.Item($row, 5).Value2 = "String"
$ws.Cells.Item($row, 6).Value2 = "String,String"

# Match the font used by the rest of the Output/Input columns in the table
# (Trebuchet MS, 10pt, black) instead of the sheet's default column style.
$fmtRange = $ws.Range("E" + $row + ":F" + $row)
$fmtRange.Font.Name = "Trebuchet MS"
$fmtRange.Font.Size = 10
$fmtRange.Font.Color = 0

$ws.Range("F23").Select() | Out-Null
